$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1) Fix the last ("position 속성") paragraph: drop the stray
#        paragraph-mark rPr (rFonts hint=eastAsia) from its pPr. ---
$last = $d.Paragraphs.Last
$lastRange = $last.Range
$fixedLastXml = "<w:p $W>" +
  "<w:pPr>" +
    "<w:pStyle w:val=""a4""/>" +
    "<w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""2""/></w:numPr>" +
    "<w:ind w:leftChars=""0""/>" +
  "</w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>p</w:t></w:r>" +
  "<w:r><w:t>osition</w:t></w:r>" +
  "<w:r><w:t xml:space=""preserve""> 속성</w:t></w:r>" +
  "</w:p>"
$lastRange.InsertXML($fixedLastXml)

# --- 2) Append the new paragraphs after it: a blank paragraph, the
#        "2022-02-24" date line, the "13강. Css 속성 3" title line, and
#        the new "float 속성" bullet paragraph. ---
$last = $d.Paragraphs.Last
$insPoint = $d.Range($last.Range.End - 1, $last.Range.End - 1)

$newXml =
  "<w:p $W/>" +
  "<w:p $W>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>2</w:t></w:r>" +
    "<w:r><w:t>022-02-24</w:t></w:r>" +
  "</w:p>" +
  "<w:p $W>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>1</w:t></w:r>" +
    "<w:r><w:t>3</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>강.</w:t></w:r>" +
    "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" +
    "<w:proofErr w:type=""spellStart""/>" +
    "<w:r><w:t>Css</w:t></w:r>" +
    "<w:proofErr w:type=""spellEnd""/>" +
    "<w:r><w:t xml:space=""preserve""> </w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t xml:space=""preserve"">속성 </w:t></w:r>" +
    "<w:r><w:t>3</w:t></w:r>" +
  "</w:p>" +
  "<w:p $W>" +
    "<w:pPr>" +
      "<w:pStyle w:val=""a4""/>" +
      "<w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""2""/></w:numPr>" +
      "<w:ind w:leftChars=""0""/>" +
      "<w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr>" +
    "</w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>f</w:t></w:r>" +
    "<w:r><w:t xml:space=""preserve"">loat </w:t></w:r>" +
    "<w:proofErr w:type=""gramStart""/>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t xml:space=""preserve"">속성 </w:t></w:r>" +
    "<w:r><w:t>:</w:t></w:r>" +
    "<w:proofErr w:type=""gramEnd""/>" +
    "<w:r><w:t xml:space=""preserve""> position</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>속성과 함께 요소의 위치를 설정하기 위한 속성.</w:t></w:r>" +
  "</w:p>"

$insPoint.InsertXML($newXml)
